# This workbook holds several embedded HCl-on-ice graph data series.
# The edit restores/fills in the "zero water" (column A) and a second
# series (column B) for rows 5-29, which previously only had columns
# C:E populated, and nudges the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2764.7464439999999, 254.11536570000001),
    @(2764.659913,        254.02469859999999),
    @(2765.025901,        254.39410430000001),
    @(2763.749433,        253.122975),
    @(2763.913787,        253.29511729999999),
    @(2764.4091469999998, 253.7955671),
    @(2764.3899799999999, 253.77290489999999),
    @(2764.140183,        253.48848520000001),
    @(2764.7834509999998, 254.1673768),
    @(2764.6191779999999, 254.00346590000001),
    @(2764.188619,        253.5656913),
    @(2764.6797499999998, 254.04858150000001),
    @(2763.5482889999998, 252.929113),
    @(2764.8275330000001, 254.24934719999999),
    @(2763.9527939999998, 253.3279167),
    @(2764.608264,        253.94857590000001),
    @(2764.2085099999999, 253.5700654),
    @(2764.188118,        253.56829569999999),
    @(2764.0114480000002, 253.3910453),
    @(2764.9836150000001, 254.3677314),
    @(2764.188619,        253.5656913),
    @(2764.1734369999999, 253.54255180000001),
    @(2765.035085,        254.40689459999999),
    @(2763.8569160000002, 253.23657969999999),
    @(2764.034506,        253.42018730000001)
)

$startRow = 5
for ($i = 0; $i -lt $values.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $values[$i][0]
    $ws.Cells.Item($r, 2).Value = $values[$i][1]
}

# Move the active-cell selection from I13 to I12, as in the edited file.
$ws.Range("I12").Select()
